# Update the "Förändrad" (changed) date column (C) for rows 2-10
# from 2023-09-11 (serial 45180) to 2023-09-12 (serial 45181).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
